# TGT.xlsx MarketBeat rank watch sheet update:
#  - Insert 3 new date columns (Jun_27, Jun_26, Jun_26) in front of the existing
#    date columns, shifting the old data from column E into column H.
#  - Add two new analyst/group rows: "Benchmark" and "Evercore ISI".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Shift the per-row "latest rating change" data (old column E) into the
#    new column H, and backfill the vacated columns E, F, G with the filler
#    value "UN" (matching every other non-event cell in the table).
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 27; $r++) {
    $oldE = $ws.Cells.Item($r, 5).Value()

    # Clear any existing highlight fill on column E before repurposing it,
    # so the highlighted rows (6 & 7) don't leave a stray green E cell behind.
    $ws.Cells.Item($r, 5).ClearFormats()

    $ws.Cells.Item($r, 5).Value = "UN"
    $ws.Cells.Item($r, 6).Value = "UN"
    $ws.Cells.Item($r, 7).Value = "UN"
    $ws.Cells.Item($r, 8).Value = $oldE
}

# Re-apply the highlight fill (light-green) to the two rows that previously
# carried it on column E -- it now belongs on column H.
$ws.Cells.Item(6, 8).Interior.ColorIndex = 35
$ws.Cells.Item(7, 8).Interior.ColorIndex = 35

# ---------------------------------------------------------------------------
# 2. Rewrite the header row. B1/C1/D1 are the three brand-new date columns;
#    E1..H1 are the previous B1..E1 headers, shifted right.
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"
$ws.Range("E1").Value = "Jun_17"
$ws.Range("F1").Value = "Jun_15"
$ws.Range("G1").Value = "Jun_13"
$ws.Range("H1").Value = "Jun_10"

# ---------------------------------------------------------------------------
# 3. Match column widths for the (now six) date columns C..H.
# ---------------------------------------------------------------------------
for ($c = 3; $c -le 8; $c++) {
    $ws.Columns.Item($c).ColumnWidth = 7.2
}

# ---------------------------------------------------------------------------
# 4. Append the two new watch-list groups at the bottom of the table.
# ---------------------------------------------------------------------------
$ws.Cells.Item(28, 1).Value = "Benchmark"
$ws.Cells.Item(28, 2).Value = "UN"
$ws.Cells.Item(28, 3).Value = "UN"
$ws.Cells.Item(28, 4).Value = "UN"

$ws.Cells.Item(29, 1).Value = "Evercore ISI"
$ws.Cells.Item(29, 2).Value = "UN"
$ws.Cells.Item(29, 3).Value = "UN"
$ws.Cells.Item(29, 4).Value = "UN"
